$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 237, shifting rows 237:255 down to 238:256.
$ws.Rows.Item(237).Insert()

# Populate the newly inserted row 237 with the new weekly record.
$ws.Range("A237").Value = 10
$ws.Range("B237").Value = "Vega Modelo de Temuco"
$ws.Range("C237").Value = "La Araucanía"
$ws.Range("D237").Value = 44461
$ws.Range("E237").Value = 9
$ws.Range("F237").Value = 100112023
$ws.Range("G237").Value = "Brócoli"
$ws.Range("H237").Value = "Sin especificar"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 2000
$ws.Range("K237").Value = 900
$ws.Range("L237").Value = 900
$ws.Range("M237").Value = 900
$ws.Range("N237").Value = "$/unidad"
$ws.Range("O237").Value = "Región de O'Higgins"
$ws.Range("P237").Value = 900
$ws.Range("Q237").Value = 1
$ws.Range("R237").Value = "Hortaliza"
